$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 347
$ws.Range("F2").Value = "11/1/2024"
$ws.Range("G2").Value = "6/1/1983"
$ws.Range("H2").Value = "1/1/2025"
$ws.Range("E3").Value = "4/1/1983"
$ws.Range("F3").Value = "11/1/2024"
$ws.Range("E4").Value = "4/1/1983"
$ws.Range("F4").Value = "11/1/2024"
$ws.Range("G4").Value = "6/1/1983"
$ws.Range("H4").Value = "1/1/2025"
$ws.Range("E5").Value = "4/1/1983"
$ws.Range("F5").Value = "11/1/2024"
$ws.Range("G5").Value = "6/1/1983"
$ws.Range("H5").Value = "1/1/2025"
$ws.Range("C6").Value = 455
$ws.Range("F6").Value = "10/1/2024"
$ws.Range("G6").Value = "6/1/1983"
$ws.Range("H6").Value = "1/1/2025"
$ws.Range("E7").Value = "3/1/1983"
$ws.Range("F7").Value = "10/1/2024"
$ws.Range("G7").Value = "6/1/1983"
$ws.Range("H7").Value = "1/1/2025"
$ws.Range("D8").Value = 423
$ws.Range("E8").Value = "3/1/1983"
$ws.Range("F8").Value = "10/1/2024"
$ws.Range("H8").Value = "1/1/2025"
$ws.Range("E9").Value = "4/1/1983"
$ws.Range("F9").Value = "11/1/2024"
$ws.Range("G9").Value = "6/1/1983"
$ws.Range("H9").Value = "1/1/2025"
$ws.Range("D10").Value = 496
$ws.Range("E10").Value = "4/1/1983"
$ws.Range("F10").Value = "11/1/2024"
$ws.Range("H10").Value = "1/1/2025"
$ws.Range("E11").Value = "3/1/1983"
$ws.Range("F11").Value = "10/1/2024"
$ws.Range("G11").Value = "6/1/1983"
$ws.Range("H11").Value = "1/1/2025"
$ws.Range("C12").Value = 384
$ws.Range("D12").Value = 365
$ws.Range("F12").Value = "11/1/2024"
$ws.Range("H12").Value = "1/1/2025"
$ws.Range("C13").Value = 480
$ws.Range("F13").Value = "11/1/2024"
$ws.Range("G13").Value = "6/1/1983"
$ws.Range("H13").Value = "1/1/2025"
$ws.Range("C14").Value = 435
$ws.Range("D14").Value = 409
$ws.Range("F14").Value = "10/1/2024"
$ws.Range("H14").Value = "1/1/2025"
$ws.Range("G15").Value = "5/2/1983"
$ws.Range("H15").Value = "1/1/2025"
$ws.Range("C16").Value = 467
$ws.Range("D16").Value = 423
$ws.Range("F16").Value = "10/1/2024"
$ws.Range("H16").Value = "1/1/2025"
$ws.Range("C17").Value = 383
$ws.Range("D17").Value = 407
$ws.Range("F17").Value = "11/1/2024"
$ws.Range("H17").Value = "1/1/2025"
$ws.Range("D18").Value = 273
$ws.Range("E18").Value = "4/1/1983"
$ws.Range("F18").Value = "11/1/2024"
$ws.Range("H18").Value = "1/1/2025"
$ws.Range("D19").Value = 411
$ws.Range("E19").Value = "4/1/1983"
$ws.Range("F19").Value = "11/1/2024"
$ws.Range("H19").Value = "1/1/2025"
$ws.Range("C20").Value = 493
$ws.Range("F20").Value = "10/1/2024"
$ws.Range("G20").Value = "6/1/1983"
$ws.Range("H20").Value = "1/1/2025"
$ws.Range("C21").Value = 323
$ws.Range("F21").Value = "11/1/2024"
$ws.Range("G21").Value = "6/1/1983"
$ws.Range("H21").Value = "1/1/2025"
$ws.Range("C22").Value = 336
$ws.Range("D22").Value = 380
$ws.Range("F22").Value = "11/1/2024"
$ws.Range("H22").Value = "1/1/2025"
$ws.Range("C23").Value = 272
$ws.Range("D23").Value = 415
$ws.Range("F23").Value = "9/1/2024"
$ws.Range("H23").Value = "1/1/2025"
$ws.Range("D24").Value = 423
$ws.Range("H24").Value = "1/1/2025"
$ws.Range("D25").Value = 328
$ws.Range("E25").Value = "3/1/1983"
$ws.Range("F25").Value = "10/1/2024"
$ws.Range("H25").Value = "1/1/2025"
$ws.Range("C26").Value = 345
$ws.Range("D26").Value = 326
$ws.Range("F26").Value = "10/1/2024"
$ws.Range("H26").Value = "1/1/2025"
$ws.Range("G27").Value = "6/1/1983"
$ws.Range("H27").Value = "1/1/2025"
$ws.Range("D28").Value = 393
$ws.Range("E28").Value = "3/1/1983"
$ws.Range("F28").Value = "10/1/2024"
$ws.Range("H28").Value = "1/1/2025"
$ws.Range("C29").Value = 274
$ws.Range("D29").Value = 249
$ws.Range("F29").Value = "10/1/2024"
$ws.Range("H29").Value = "1/1/2025"
$ws.Range("D30").Value = 230
$ws.Range("E30").Value = "4/1/1983"
$ws.Range("F30").Value = "11/1/2024"
$ws.Range("H30").Value = "1/1/2025"
$ws.Range("C31").Value = 406
$ws.Range("F31").Value = "10/1/2024"
$ws.Range("G31").Value = "6/1/1983"
$ws.Range("H31").Value = "1/1/2025"
$ws.Range("E32").Value = "4/1/1983"
$ws.Range("F32").Value = "11/1/2024"
$ws.Range("G32").Value = "6/1/1983"
$ws.Range("H32").Value = "1/1/2025"
$ws.Range("D33").Value = 411
$ws.Range("H33").Value = "1/1/2025"
$ws.Range("C34").Value = 214
$ws.Range("D34").Value = 330
$ws.Range("F34").Value = "10/1/2024"
$ws.Range("H34").Value = "1/1/2025"
$ws.Range("C35").Value = 417
$ws.Range("D35").Value = 330
$ws.Range("F35").Value = "10/1/2024"
$ws.Range("H35").Value = "1/1/2025"
$ws.Range("D36").Value = 423
$ws.Range("E36").Value = "4/1/1983"
$ws.Range("F36").Value = "11/1/2024"
$ws.Range("H36").Value = "1/1/2025"
$ws.Range("C37").Value = 479
$ws.Range("D37").Value = 330
$ws.Range("F37").Value = "11/1/2024"
$ws.Range("H37").Value = "1/1/2025"
$ws.Range("C38").Value = 372
$ws.Range("D38").Value = 380
$ws.Range("F38").Value = "11/1/2024"
$ws.Range("H38").Value = "1/1/2025"
$ws.Range("C39").Value = 240
$ws.Range("D39").Value = 236
$ws.Range("F39").Value = "11/1/2024"
$ws.Range("H39").Value = "1/1/2025"
$ws.Range("C40").Value = 299
$ws.Range("D40").Value = 327
$ws.Range("F40").Value = "11/1/2024"
$ws.Range("H40").Value = "1/1/2025"
$ws.Range("C41").Value = 406
$ws.Range("D41").Value = 249
$ws.Range("F41").Value = "10/1/2024"
$ws.Range("H41").Value = "1/1/2025"
$ws.Range("C42").Value = 249
$ws.Range("D42").Value = 231
$ws.Range("F42").Value = "10/1/2024"
$ws.Range("H42").Value = "1/1/2025"
$ws.Range("C43").Value = 485
$ws.Range("D43").Value = 330
$ws.Range("F43").Value = "10/1/2024"
$ws.Range("H43").Value = "1/1/2025"
$ws.Range("C44").Value = 418
$ws.Range("D44").Value = 317
$ws.Range("F44").Value = "10/1/2024"
$ws.Range("H44").Value = "1/1/2025"
$ws.Range("D45").Value = 330
$ws.Range("H45").Value = "1/1/2025"
$ws.Range("C46").Value = 346
$ws.Range("D46").Value = 311
$ws.Range("F46").Value = "10/1/2024"
$ws.Range("H46").Value = "1/1/2025"
$ws.Range("C47").Value = 348
$ws.Range("D47").Value = 273
$ws.Range("F47").Value = "11/1/2024"
$ws.Range("H47").Value = "1/1/2025"
$ws.Range("C48").Value = 371
$ws.Range("D48").Value = 328
$ws.Range("F48").Value = "11/1/2024"
$ws.Range("H48").Value = "1/1/2025"
$ws.Range("C49").Value = 310
$ws.Range("D49").Value = 326
$ws.Range("F49").Value = "10/1/2024"
$ws.Range("H49").Value = "1/1/2025"
$ws.Range("C50").Value = 371
$ws.Range("D50").Value = 251
$ws.Range("F50").Value = "11/1/2024"
$ws.Range("H50").Value = "1/1/2025"
$ws.Range("D51").Value = 330
$ws.Range("H51").Value = "1/1/2025"
$ws.Range("C52").Value = 354
$ws.Range("D52").Value = 328
$ws.Range("F52").Value = "10/1/2024"
$ws.Range("H52").Value = "1/1/2025"
